$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append the new "2023" data row (row 24) under the existing economic-data
# table, following the same layout/styling pattern as the prior rows.
# ---------------------------------------------------------------------------

# 1) Seed row 24 from row 23 so it inherits the same number formats, fonts,
#    and borders as the rest of the data table (reuses the existing styles
#    instead of inventing mismatched ones).
$ws.Range("A23:T23").Copy($ws.Range("A24:T24"))

# 2) Write the new values for 2023.
$ws.Range("A24").Value2 = 2023
$ws.Range("B24").Value2 = -0.57999999999999996
$ws.Range("C24").Value2 = -1.93
$ws.Range("D24").Value2 = -0.99
$ws.Range("E24").Value2 = 3
$ws.Range("F24").Value2 = -0.86
$ws.Range("G24").ClearContents()
$ws.Range("H24").Value2 = 5.4080000000000004
$ws.Range("I24").Value2 = 1616.3068878659999
$ws.Range("J24").Value2 = 1365.2774322242999
$ws.Range("K24").Value2 = -0.0396925228069165
$ws.Range("L24").Value2 = 40953671000000
$ws.Range("M24").Value2 = 30.768128065700498
$ws.Range("N24").Value2 = -0.24421476876381301
$ws.Range("O24").Value2 = -825231745.482939
$ws.Range("P24").Value2 = 43.260843757783597
$ws.Range("Q24").Value2 = 42.571738430231598
$ws.Range("R24").Value2 = 16.253063175299999
$ws.Range("S24").Value2 = 280.35611165514001
$ws.Range("T24").Value2 = 91.394551203056807

# 3) A24 (the year label) gets the header-style look used elsewhere in the
#    sheet: bold Cambria, thin border, horizontally centered only (no
#    vertical centering, no wrap) -- matches column A's other year cells in
#    weight/face but with single-axis centering.
$ws.Range("B2").Copy($ws.Range("A24"))
$ws.Range("A24").Value2 = 2023
$ws.Range("A24").WrapText = $false
$ws.Range("A24").VerticalAlignment = -4107

# 4) O24 uses a 2-decimal numeric format (like most of the row) rather than
#    the integer format O-column sometimes uses.
$ws.Range("O24").NumberFormat = "0.00"

# 5) G24 has no data for 2023; give it a plain (unbolded, unformatted)
#    bordered cell instead of the numeric style copied from row 23.
$ws.Range("G24").ClearFormats()
$ws.Range("G24").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# View-state touch-ups (zoom level / selected cell) to reflect the state the
# workbook was left in.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 55
$ws.Range("I27").Select()
